$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1149.9333
$ws.Range("J17").Value = 1165.841
$ws.Range("L17").Value = 3497.523
$ws.Range("N17").Value = -3833.523
$ws.Range("H53").Value = 606.8077
$ws.Range("I53").Value = 318.23077
$ws.Range("J53").Value = 895.38464
$ws.Range("K53").Value = 318.23077
$ws.Range("L53").Value = 895.38464
$ws.Range("M53").Value = 318.76923
$ws.Range("N53").Value = -2169.38464
$ws.Range("H116").Value = 1600
$ws.Range("I116").Value = 1466.6666
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1466.6666
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1975.3334
$ws.Range("N116").Value = -8884
$ws.Range("H132").Value = 2370.0862
$ws.Range("I132").Value = 2738.9148
$ws.Range("J132").Value = 794.1818
$ws.Range("K132").Value = 8216.7444
$ws.Range("L132").Value = 2382.5454
$ws.Range("M132").Value = -5686.7444
$ws.Range("N132").Value = -7442.5454
$ws.Range("H137").Value = 14706897
$ws.Range("I137").Value = 911.3617
$ws.Range("K137").Value = 2734.0851
$ws.Range("M137").Value = -184.0851000000002
$ws.Range("H138").Value = 2353.6267
$ws.Range("I138").Value = 1743.2195
$ws.Range("J138").Value = 3089.7058
$ws.Range("K138").Value = 5229.6585
$ws.Range("L138").Value = 9269.117400000001
$ws.Range("M138").Value = -89.65849999999955
$ws.Range("N138").Value = -19549.1174

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1588195.8
$ws.Range("I61").Value = 1822381.9
$ws.Range("J61").Value = 933.55554
$ws.Range("K61").Value = 1822381.9
$ws.Range("L61").Value = 933.55554
$ws.Range("M61").Value = -1822169.9
$ws.Range("N61").Value = -1357.55554
$ws.Range("H80").Value = 20000.666
$ws.Range("J80").Value = 20000.666
$ws.Range("L80").Value = 20000.666
$ws.Range("N80").Value = -21996.666
$ws.Range("H83").Value = 20000.666
$ws.Range("J83").Value = 20000.666
$ws.Range("L83").Value = 60001.99800000001
$ws.Range("N83").Value = -69985.99800000001
$ws.Range("H132").Value = 4924024.5
$ws.Range("I132").Value = 5533936
$ws.Range("K132").Value = 16601808
$ws.Range("M132").Value = -16599278
$ws.Range("H136").Value = 1588195.8
$ws.Range("I136").Value = 1822381.9
$ws.Range("J136").Value = 933.55554
$ws.Range("K136").Value = 5467145.699999999
$ws.Range("L136").Value = 2800.66662
$ws.Range("M136").Value = -5464595.699999999
$ws.Range("N136").Value = -7900.66662

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2833147.2
$ws.Range("I134").Value = 3657074.2
$ws.Range("J134").Value = 899
$ws.Range("K134").Value = 10971222.6
$ws.Range("L134").Value = 2697
$ws.Range("M134").Value = -10968687.6
$ws.Range("N134").Value = -7767

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2637835.2
$ws.Range("I31").Value = 944.5472
$ws.Range("J31").Value = 10858730
$ws.Range("K31").Value = 944.5472
$ws.Range("L31").Value = 10858730
$ws.Range("M31").Value = -649.5472
$ws.Range("N31").Value = -10859320
$ws.Range("H34").Value = 2637835.2
$ws.Range("I34").Value = 944.5472
$ws.Range("J34").Value = 10858730
$ws.Range("K34").Value = 944.5472
$ws.Range("L34").Value = 10858730
$ws.Range("M34").Value = -742.5472
$ws.Range("N34").Value = -10859134
$ws.Range("H58").Value = 1508.4124
$ws.Range("I58").Value = 822.5323
$ws.Range("J58").Value = 2723.4
$ws.Range("K58").Value = 822.5323
$ws.Range("L58").Value = 2723.4
$ws.Range("M58").Value = -619.5323
$ws.Range("N58").Value = -3129.4
$ws.Range("H132").Value = 1461.9365
$ws.Range("I132").Value = 1445.0944
$ws.Range("J132").Value = 1551.2
$ws.Range("K132").Value = 4335.2832
$ws.Range("L132").Value = 4653.6
$ws.Range("M132").Value = -1805.2832
$ws.Range("N132").Value = -9713.6
$ws.Range("H134").Value = 1225.8654
$ws.Range("I134").Value = 1274.2
$ws.Range("J134").Value = 915.1429000000001
$ws.Range("K134").Value = 3822.6
$ws.Range("L134").Value = 2745.4287
$ws.Range("M134").Value = -1287.6
$ws.Range("N134").Value = -7815.4287
$ws.Range("H136").Value = 1508.4124
$ws.Range("I136").Value = 822.5323
$ws.Range("J136").Value = 2723.4
$ws.Range("K136").Value = 2467.5969
$ws.Range("L136").Value = 8170.200000000001
$ws.Range("M136").Value = 82.40309999999999
$ws.Range("N136").Value = -13270.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10205266
$ws.Range("I122").Value = 20834010
$ws.Range("J122").Value = 1254745.6
$ws.Range("K122").Value = 187506090
$ws.Range("L122").Value = 11292710.4
$ws.Range("M122").Value = -187503640
$ws.Range("N122").Value = -11297610.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4322
$ws.Range("I70").Value = 4239.8
$ws.Range("J70").Value = 4424.75
$ws.Range("K70").Value = 4239.8
$ws.Range("L70").Value = 4424.75
$ws.Range("M70").Value = -3969.8
$ws.Range("N70").Value = -4964.75
$ws.Range("H73").Value = 4322
$ws.Range("I73").Value = 4239.8
$ws.Range("J73").Value = 4424.75
$ws.Range("K73").Value = 4239.8
$ws.Range("L73").Value = 4424.75
$ws.Range("M73").Value = -3303.8
$ws.Range("N73").Value = -6296.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3475.4443
$ws.Range("I62").Value = 3599.6667
$ws.Range("J62").Value = 3413.3333
$ws.Range("K62").Value = 3599.6667
$ws.Range("L62").Value = 3413.3333
$ws.Range("M62").Value = -2975.6667
$ws.Range("N62").Value = -4661.3333
$ws.Range("H65").Value = 3475.4443
$ws.Range("I65").Value = 3599.6667
$ws.Range("J65").Value = 3413.3333
$ws.Range("K65").Value = 17998.3335
$ws.Range("L65").Value = 17066.6665
$ws.Range("M65").Value = -14878.3335
$ws.Range("N65").Value = -23306.6665
$ws.Range("H96").Value = 5566769.5
$ws.Range("I96").Value = 14287970
$ws.Range("J96").Value = 16914.727
$ws.Range("K96").Value = 14287970
$ws.Range("L96").Value = 16914.727
$ws.Range("M96").Value = -14286597
$ws.Range("N96").Value = -19660.727
$ws.Range("H107").Value = 608.9091
$ws.Range("I107").Value = 226
$ws.Range("J107").Value = 827.7143
$ws.Range("K107").Value = 678
$ws.Range("L107").Value = 2483.1429
$ws.Range("M107").Value = 1242
$ws.Range("N107").Value = -6323.1429
$ws.Range("H113").Value = 453.5
$ws.Range("I113").Value = 413.7143
$ws.Range("J113").Value = 523.125
$ws.Range("K113").Value = 1241.1429
$ws.Range("L113").Value = 1569.375
$ws.Range("M113").Value = 928.8571000000002
$ws.Range("N113").Value = -5909.375
